$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking "Price" text values keep their original text formatting
# (all Price/Volume columns are text in the source data, e.g. "379.10", "0.580").

$ws.Range('D2').Value = '50.934.79'
$ws.Range('E2').Value = '  -0.51%  '
$ws.Range('D3').Value = '2.945.91'
$ws.Range('E3').Value = '  -0.22%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '379.10'
$ws.Range('E5').Value = '  -0.03%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '100.86'
$ws.Range('E6').Value = '  -1.72%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.539'
$ws.Range('E7').Value = '  -0.44%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.580'
$ws.Range('E9').Value = '  -1.65%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '36.11'
$ws.Range('E10').Value = '  -1.19%  '
$ws.Range('E11').Value = '  -0.44%  '
$ws.Range('E12').Value = '  +1.53%  '
$ws.Range('D13').Value = '3.396.17'
$ws.Range('E13').Value = '  -0.78%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '18.22'
$ws.Range('E14').Value = '  +1.29%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.62'
$ws.Range('E15').Value = '  +3.02%  '
$ws.Range('B16').Value = 'Uniswap'
$ws.Range('C16').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '12.07'
$ws.Range('E16').Value = '  +68.24%  '
$ws.Range('D17').Value = '2.943.25'
$ws.Range('E17').Value = '  -0.25%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.995'
$ws.Range('E18').Value = '  +1.06%  '
$ws.Range('D19').Value = '50.901.78'
$ws.Range('E19').Value = '  -0.73%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '3.03'
$ws.Range('E20').Value = '  -5.66%  '
$ws.Range('E21').Value = '  -1.42%  '
$ws.Range('D22').Value = '0.0₃0949'
$ws.Range('E22').Value = '  -0.29%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '69.51'
$ws.Range('E23').Value = '  +1.70%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '266.87'
$ws.Range('E24').Value = '  +2.05%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.26'
$ws.Range('E25').Value = '  +13.73%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '8.07'
$ws.Range('E26').Value = '  -2.65%  '
$ws.Range('E27').Value = '  -0.05%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '7.10'
$ws.Range('E28').Value = '  -6.55%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '25.60'
$ws.Range('E29').Value = '  -0.39%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.162'
$ws.Range('E30').Value = '  -3.73%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.109'
$ws.Range('E31').Value = '  -2.88%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '10.02'
$ws.Range('E32').Value = '  +2.41%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '50.46'
$ws.Range('E33').Value = '  -0.05%  '
$ws.Range('E34').Value = '  +0.33%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '33.51'
$ws.Range('E35').Value = '  -0.87%  '
$ws.Range('E36').Value = '  -3.07%  '
$ws.Range('E37').Value = '  -0.17%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.11'
$ws.Range('E38').Value = '  +4.43%  '
$ws.Range('E39').Value = '  +1.19%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '16.57'
$ws.Range('E40').Value = '  -2.04%  '
$ws.Range('E41').Value = '  +2.16%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.50'
$ws.Range('E42').Value = '  -2.29%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '119.96'
$ws.Range('E43').Value = '  -1.48%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '21.26'
$ws.Range('E44').Value = '  +0.49%  '
$ws.Range('E45').Value = '  +6.83%  '
$ws.Range('E46').Value = '  -1.84%  '
$ws.Range('E47').Value = '  -1.20%  '
$ws.Range('D48').Value = '2.010.78'
$ws.Range('E48').Value = '  +0.46%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.260'
$ws.Range('E49').Value = '  -4.09%  '
$ws.Range('E50').Value = '  -6.25%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '5.30'
$ws.Range('E51').Value = '  +4.67%  '

Write-Host "Updated cryptos list"
